$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-looking-like-numbers to stay as text (matches inlineStr in source) by
# pre-formatting the affected cells as Text before assigning their values.
$textCells = @(
    "D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "E48", "D49", "E49", "D50", "E50", "D51", "E51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "304.58"
$ws.Range("E2").Value = "0.13%"
$ws.Range("D3").Value = "35.51"
$ws.Range("E3").Value = "-4.39%"
$ws.Range("E4").Value = "0.29%"
$ws.Range("D5").Value = "0.07862"
$ws.Range("E5").Value = "0.14%"
$ws.Range("D6").Value = "2.131"
$ws.Range("E6").Value = "-3.49%"
$ws.Range("D7").Value = "7.939"
$ws.Range("E7").Value = "-0.93%"
$ws.Range("D8").Value = "0.9251"
$ws.Range("E8").Value = "0.12%"
$ws.Range("D9").Value = "0.09741"
$ws.Range("E9").Value = "-1.19%"
$ws.Range("D10").Value = "0.1837"
$ws.Range("E10").Value = "-2.77%"
$ws.Range("D11").Value = "0.08615"
$ws.Range("E11").Value = "-0.32%"
$ws.Range("D12").Value = "0.03587"
$ws.Range("E12").Value = "-2.35%"
$ws.Range("D13").Value = "0.09942"
$ws.Range("E13").Value = "0.02%"
$ws.Range("D14").Value = "0.001447"
$ws.Range("E14").Value = "-2.20%"
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D15").Value = "0.04566"
$ws.Range("E15").Value = "-0.77%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.005727"
$ws.Range("E16").Value = "1.53%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.474"
$ws.Range("E17").Value = "0.37%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "4.140"
$ws.Range("E18").Value = "3.00%"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "2.752"
$ws.Range("E19").Value = "22.11%"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "0.3372"
$ws.Range("E20").Value = "-1.14%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "0.1349"
$ws.Range("E21").Value = "3.58%"
$ws.Range("B22").Value = "MCDex"
$ws.Range("C22").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D22").Value = "5.163"
$ws.Range("E22").Value = "8.27%"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").Value = "0.2213"
$ws.Range("E23").Value = "0.25%"
$ws.Range("D24").Value = "0.001233"
$ws.Range("E24").Value = "-1.88%"
$ws.Range("D25").Value = "0.004825"
$ws.Range("E25").Value = "7.62%"
$ws.Range("D26").Value = "0.0001300"
$ws.Range("E26").Value = "-7.27%"
$ws.Range("D27").Value = "0.0004761"
$ws.Range("E27").Value = "74.98%"
$ws.Range("D39").Value = "0.01847"
$ws.Range("E39").Value = "0.23%"
$ws.Range("D40").Value = "0.04729"
$ws.Range("E40").Value = "-0.76%"
$ws.Range("D41").Value = "0.007788"
$ws.Range("E41").Value = "-3.53%"
$ws.Range("D42").Value = "0.1386"
$ws.Range("E42").Value = "-1.07%"
$ws.Range("D43").Value = "0.007765"
$ws.Range("D44").Value = "0.002163"
$ws.Range("E44").Value = "2.97%"
$ws.Range("D45").Value = "0.01133"
$ws.Range("E45").Value = "8.92%"
$ws.Range("D46").Value = "0.00006301"
$ws.Range("E46").Value = "0.08%"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").Value = "-0.12%"
$ws.Range("E48").Value = "-0.03%"
$ws.Range("D49").Value = "50.71"
$ws.Range("E49").Value = "32.71%"
$ws.Range("D50").Value = "0.001905"
$ws.Range("E50").Value = "-29.26%"
$ws.Range("D51").Value = "0.00002101"
$ws.Range("E51").Value = "-0.12%"
